# "download articles with pandoc title blocks"
#
# Turns the old two-paragraph title block:
#   Heading1 "On Pilgrimage - May 1961"
#   Bold     "By Dorothy Day"
# into a pandoc-style title block:
#   Title   "On" " " "Pilgrimage" " " "-" " " "May" " " "1961"   (Title style, one run per word/space)
#   Authors "Dorothy" " " "Day"                                  (Authors style, "By " dropped)

$d = $word.ActiveDocument

function New-RunsXml($words) {
    $xml = ""
    foreach ($w in $words) {
        $xml += '<w:r><w:t xml:space="preserve">' + $w + '</w:t></w:r>'
    }
    return $xml
}

function New-PackageXml($bodyInner) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $bodyInner + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

$titleRuns   = New-RunsXml @("On", " ", "Pilgrimage", " ", "-", " ", "May", " ", "1961")
$authorRuns  = New-RunsXml @("Dorothy", " ", "Day")

$titleParaXml   = '<w:p><w:pPr><w:pStyle w:val="Title"/></w:pPr>' + $titleRuns + '</w:p>'
$authorsParaXml = '<w:p><w:pPr><w:pStyle w:val="Authors"/></w:pPr>' + $authorRuns + '</w:p>'

# The document currently starts with:
#   Paragraph 1 (Heading1): "On Pilgrimage - May 1961"
#   Paragraph 2 (bold run): "By Dorothy Day"
# InsertXML replaces exactly the contents of the range it's called on, so
# target each paragraph's Range individually - this also keeps any
# surrounding (non-editable-via-this-API) structural markers, such as the
# bookmark wrapping paragraph 1, scoped to that same paragraph rather than
# spreading across both, same as in the original document.
# Update paragraph 2 first so paragraph 1's replacement doesn't shift its index/offsets.
$secondPara = $d.Paragraphs(2)
[void]($secondPara.Range.InsertXML((New-PackageXml $authorsParaXml)))

$firstPara = $d.Paragraphs(1)
[void]($firstPara.Range.InsertXML((New-PackageXml $titleParaXml)))
